$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the date-only number format currently used by the last row (Y28),
# since the new row 29 will become the new last row and should use it.
$dateOnlyFormat = $ws.Range("Y28").NumberFormat

# Row 28 is no longer the last row of data -> its Date cell switches from the
# date-only format to the regular date-time format used by all other rows.
$ws.Range("Y28").NumberFormat = $ws.Range("Y27").NumberFormat

# Append the new row of bunker price data (row 29)
$row29 = @(565,474,450,522,515,520,474,570,490,450,571,480,485,505,545,480,618,490,474,480,619,550,599,495,45754,850,555,543.5,500,545,507,509,745,473,735,474,488,570,555,488,535,547,568,547,645,634,496,485)

for ($i = 0; $i -lt $row29.Length; $i++) {
    $ws.Cells.Item(29, $i + 1).Value = $row29[$i]
}

# The new last row's Date cell (Y29) gets the date-only format
$ws.Range("Y29").NumberFormat = $dateOnlyFormat
